$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D2:E51 to Text format up front so numeric-looking strings
# (e.g. "29.881.48", "1.005", "0.000007740") are stored verbatim
# instead of being auto-converted to numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

# Row 2
$ws.Range('D2').Value = '29.881.48'
$ws.Range('E2').Value = '  +0.07%  '

# Row 3
$ws.Range('D3').Value = '1.899.13'
$ws.Range('E3').Value = '  -0.17%  '

# Row 5
$ws.Range('D5').Value = '0.7933'
$ws.Range('E5').Value = '  -0.90%  '

# Row 6
$ws.Range('D6').Value = '242.95'
$ws.Range('E6').Value = '  +0.93%  '

# Row 7
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.51%  '

# Row 8
$ws.Range('D8').Value = '0.3203'
$ws.Range('E8').Value = '  +2.61%  '

# Row 9
$ws.Range('D9').Value = '26.23'
$ws.Range('E9').Value = '  +0.62%  '

# Row 10
$ws.Range('D10').Value = '0.07106'
$ws.Range('E10').Value = '  +3.47%  '

# Row 11
$ws.Range('D11').Value = '0.08057'
$ws.Range('E11').Value = '  +1.07%  '

# Row 12
$ws.Range('D12').Value = '0.7707'
$ws.Range('E12').Value = '  +5.17%  '

# Row 13
$ws.Range('D13').Value = '1.959.19'
$ws.Range('E13').Value = '  +2.86%  '

# Row 14
$ws.Range('D14').Value = '5.315'
$ws.Range('E14').Value = '  +3.14%  '

# Row 15
$ws.Range('D15').Value = '92.59'
$ws.Range('E15').Value = '  +0.17%  '

# Row 16
$ws.Range('D16').Value = '29.929.94'
$ws.Range('E16').Value = '  +0.17%  '

# Row 17
$ws.Range('D17').Value = '13.89'
$ws.Range('E17').Value = '  +0.18%  '

# Row 18
$ws.Range('D18').Value = '5.921'
$ws.Range('E18').Value = '  +1.44%  '

# Row 19
$ws.Range('D19').Value = '245.05'
$ws.Range('E19').Value = '  +0.42%  '

# Row 20
$ws.Range('D20').Value = '0.000007740'
$ws.Range('E20').Value = '  +0.72%  '

# Row 21
$ws.Range('D21').Value = '2.167.21'
$ws.Range('E21').Value = '  +0.33%  '

# Row 22
$ws.Range('D22').Value = '1.005'
$ws.Range('E22').Value = '  +0.46%  '

# Row 23
$ws.Range('D23').Value = '8.144'
$ws.Range('E23').Value = '  +18.61%  '

# Row 24
$ws.Range('D24').Value = '1.007'
$ws.Range('E24').Value = '  +0.49%  '

# Row 25
$ws.Range('D25').Value = '0.1622'
$ws.Range('E25').Value = '  +15.04%  '

# Row 26
$ws.Range('D26').Value = '9.313'
$ws.Range('E26').Value = '  +1.62%  '

# Row 27
$ws.Range('D27').Value = '165.74'
$ws.Range('E27').Value = '  -0.84%  '

# Row 28
$ws.Range('D28').Value = '18.73'
$ws.Range('E28').Value = '  -0.34%  '

# Row 29
$ws.Range('D29').Value = '2.094'
$ws.Range('E29').Value = '  +4.20%  '

# Row 30
$ws.Range('D30').Value = '1.377'
$ws.Range('E30').Value = '  +1.36%  '

# Row 31
$ws.Range('D31').Value = '1.539'
$ws.Range('E31').Value = '  +1.61%  '

# Row 32
$ws.Range('D32').Value = '4.481'
$ws.Range('E32').Value = '  +4.98%  '

# Row 33
$ws.Range('D33').Value = '0.05685'
$ws.Range('E33').Value = '  +2.81%  '

# Row 34
$ws.Range('D34').Value = '4.082'
$ws.Range('E34').Value = '  +0.84%  '

# Row 35
$ws.Range('D35').Value = '1.265'
$ws.Range('E35').Value = '  +1.09%  '

# Row 36
$ws.Range('D36').Value = '0.7367'
$ws.Range('E36').Value = '  +1.29%  '

# Row 37
$ws.Range('D37').Value = '1.002'
$ws.Range('E37').Value = '  +0.21%  '

# Row 38
$ws.Range('D38').Value = '2.718'
$ws.Range('E38').Value = '  -0.19%  '

# Row 39
$ws.Range('D39').Value = '0.01931'
$ws.Range('E39').Value = '  +0.76%  '

# Row 40
$ws.Range('D40').Value = '2.783'
$ws.Range('E40').Value = '  -0.06%  '

# Row 41
$ws.Range('E41').Value = '  +1.69%  '

# Row 42
$ws.Range('D42').Value = '72.29'
$ws.Range('E42').Value = '  +0.82%  '

# Row 43
$ws.Range('D43').Value = '5.926'
$ws.Range('E43').Value = '  -1.62%  '

# Row 44
$ws.Range('D44').Value = '0.8469'
$ws.Range('E44').Value = '  +1.57%  '

# Row 45
$ws.Range('D45').Value = '1.005'
$ws.Range('E45').Value = '  +0.53%  '

# Row 46
$ws.Range('D46').Value = '1.888'
$ws.Range('E46').Value = '  +1.49%  '

# Row 47
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '102.27'
$ws.Range('E47').Value = '  +1.94%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.022.71'
$ws.Range('E48').Value = '  +4.92%  '

# Row 49
$ws.Range('D49').Value = '9.823'
$ws.Range('E49').Value = '  +1.73%  '

# Row 50
$ws.Range('D50').Value = '7.496'
$ws.Range('E50').Value = '  -0.40%  '

# Row 51
$ws.Range('D51').Value = '3.027'
$ws.Range('E51').Value = '  +9.54%  '
